$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 37, shifting rows 37:151 down to 38:152
$ws.Rows("37:37").Insert()

# Populate the newly inserted row 37 with the required data.
# Most values are the same as the old row 37 had, except D (Fecha), K/L/M (prices) and P (Precio $/Kg).
$ws.Range("A37").Value2 = 5
$ws.Range("B37").Value2 = 'Macroferia Regional de Talca'
$ws.Range("C37").Value2 = 'Maule'
$ws.Range("D37").Value2 = 44998
$ws.Range("E37").Value2 = 7
$ws.Range("F37").Value2 = 100112001
$ws.Range("G37").Value2 = 'Berenjena'
$ws.Range("H37").Value2 = 'Sin especificar'
$ws.Range("I37").Value2 = 'Primera'
$ws.Range("J37").Value2 = 200
$ws.Range("K37").Value2 = 9000
$ws.Range("L37").Value2 = 9000
$ws.Range("M37").Value2 = 9000
$ws.Range("N37").Value2 = '$/caja 60 unidades'
$ws.Range("O37").Value2 = 'Región del Maule'
$ws.Range("P37").Value2 = 150
$ws.Range("Q37").Value2 = 60
$ws.Range("R37").Value2 = 'Hortaliza'

# Ensure the date cell keeps the same date/time number format as the rest of column D
$ws.Range("D37").NumberFormat = $ws.Range("D36").NumberFormat
